$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H53").Value = 3745.4119
$ws.Range("I53").Value = 2318.1428
$ws.Range("K53").Value = 2318.1428
$ws.Range("M53").Value = -1681.1428
$ws.Range("H92").Value = 654.4545000000001
$ws.Range("I92").Value = 600
$ws.Range("J92").Value = 899.5
$ws.Range("K92").Value = 600
$ws.Range("L92").Value = 899.5
$ws.Range("M92").Value = 648
$ws.Range("N92").Value = -3395.5
$ws.Range("H132").Value = 1290.2916
$ws.Range("I132").Value = 1390.1578
$ws.Range("J132").Value = 910.8
$ws.Range("K132").Value = 4170.4734
$ws.Range("L132").Value = 2732.4
$ws.Range("M132").Value = -1640.4734
$ws.Range("N132").Value = -7792.4
$ws.Range("H135").Value = 278303.47
$ws.Range("I135").Value = 286112.16
$ws.Range("J135").Value = 5000
$ws.Range("K135").Value = 2575009.44
$ws.Range("L135").Value = 45000
$ws.Range("M135").Value = -2572474.44
$ws.Range("N135").Value = -50070

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H10").Value = 5500
$ws.Range("I10").Value = 5500
$ws.Range("K10").Value = 5500
$ws.Range("M10").Value = -5330
$ws.Range("H16").Value = 1964.8334
$ws.Range("J16").Value = 2833
$ws.Range("L16").Value = 2833
$ws.Range("N16").Value = -3407
$ws.Range("H44").Value = 30000
$ws.Range("J44").Value = 30000
$ws.Range("L44").Value = 30000
$ws.Range("N44").Value = -30976
$ws.Range("H55").Value = 30000
$ws.Range("I55").Value = 0
$ws.Range("J55").Value = 30000
$ws.Range("K55").Value = 0
$ws.Range("L55").Value = 30000
$ws.Range("M55").ClearContents()
$ws.Range("N55").Value = -30630
$ws.Range("H74").Value = 26560.742
$ws.Range("I74").Value = 33979.367
$ws.Range("J74").Value = 5177.647
$ws.Range("K74").Value = 33979.367
$ws.Range("L74").Value = 5177.647
$ws.Range("M74").Value = -33105.367
$ws.Range("N74").Value = -6925.647
$ws.Range("H77").Value = 26560.742
$ws.Range("I77").Value = 33979.367
$ws.Range("J77").Value = 5177.647
$ws.Range("K77").Value = 169896.835
$ws.Range("L77").Value = 25888.235
$ws.Range("M77").Value = -165528.835
$ws.Range("N77").Value = -34624.235
$ws.Range("H122").Value = 18626.467
$ws.Range("I122").Value = 24239.9
$ws.Range("K122").Value = 72719.70000000001
$ws.Range("M122").Value = -70269.70000000001

# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H24").Value = 0
$ws.Range("I24").Value = 0
$ws.Range("K24").Value = 0
$ws.Range("M24").ClearContents()
$ws.Range("H25").Value = 1688.3334
$ws.Range("I25").Value = 1688.3334
$ws.Range("K25").Value = 1688.3334
$ws.Range("M25").Value = -1453.3334
$ws.Range("H29").Value = 50449
$ws.Range("I29").Value = 899
$ws.Range("K29").Value = 899
$ws.Range("M29").Value = -610
$ws.Range("H34").Value = 0
$ws.Range("I34").Value = 0
$ws.Range("K34").Value = 0
$ws.Range("M34").ClearContents()
$ws.Range("H37").Value = 6562.3335
$ws.Range("I37").Value = 0
$ws.Range("J37").Value = 6562.3335
$ws.Range("K37").Value = 0
$ws.Range("L37").Value = 6562.3335
$ws.Range("M37").ClearContents()
$ws.Range("N37").Value = -6836.3335
$ws.Range("H39").Value = 18000
$ws.Range("J39").Value = 18000
$ws.Range("L39").Value = 18000
$ws.Range("N39").Value = -18778
$ws.Range("H99").Value = 1814.4
$ws.Range("I99").Value = 790.6667
$ws.Range("K99").Value = 790.6667
$ws.Range("M99").Value = 707.3333
$ws.Range("H107").Value = 83337050
$ws.Range("I107").Value = 125002710
$ws.Range("J107").Value = 5718.75
$ws.Range("K107").Value = 125002710
$ws.Range("L107").Value = 5718.75
$ws.Range("M107").Value = -125000790
$ws.Range("N107").Value = -9558.75
$ws.Range("H134").Value = 3643.6353
$ws.Range("I134").Value = 2683.6035
$ws.Range("J134").Value = 7123.75
$ws.Range("K134").Value = 8050.810500000001
$ws.Range("L134").Value = 21371.25
$ws.Range("M134").Value = -5515.810500000001
$ws.Range("N134").Value = -26441.25

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H33").Value = 7517.5
$ws.Range("I33").Value = 5000
$ws.Range("K33").Value = 5000
$ws.Range("M33").Value = -4621
$ws.Range("H35").Value = 30028
$ws.Range("I35").Value = 0
$ws.Range("J35").Value = 30028
$ws.Range("K35").Value = 0
$ws.Range("L35").Value = 30028
$ws.Range("M35").ClearContents()
$ws.Range("N35").Value = -30616
$ws.Range("H36").Value = 35000
$ws.Range("J36").Value = 50000
$ws.Range("L36").Value = 50000
$ws.Range("N36").Value = -50776
$ws.Range("H37").Value = 0
$ws.Range("J37").Value = 0
$ws.Range("L37").Value = 0
$ws.Range("N37").ClearContents()
$ws.Range("H39").Value = 26000
$ws.Range("I39").Value = 26000
$ws.Range("J39").Value = 0
$ws.Range("K39").Value = 26000
$ws.Range("L39").Value = 0
$ws.Range("M39").Value = -25609
$ws.Range("N39").ClearContents()
$ws.Range("H40").Value = 35000
$ws.Range("J40").Value = 50000
$ws.Range("L40").Value = 50000
$ws.Range("N40").Value = -50320
$ws.Range("H49").Value = 26000
$ws.Range("I49").Value = 26000
$ws.Range("J49").Value = 0
$ws.Range("K49").Value = 26000
$ws.Range("L49").Value = 0
$ws.Range("M49").Value = -25818
$ws.Range("N49").ClearContents()
$ws.Range("H59").Value = 98250
$ws.Range("J59").Value = 98250
$ws.Range("L59").Value = 98250
$ws.Range("N59").Value = -100540
$ws.Range("H132").Value = 4395.6665
$ws.Range("I132").Value = 1942.9286
$ws.Range("K132").Value = 5828.7858
$ws.Range("M132").Value = -3298.7858
$ws.Range("H134").Value = 6772.9795
$ws.Range("I134").Value = 7369.7407
$ws.Range("J134").Value = 6040.591
$ws.Range("K134").Value = 22109.2221
$ws.Range("L134").Value = 18121.773
$ws.Range("M134").Value = -19574.2221
$ws.Range("N134").Value = -23191.773

# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H137").Value = 251228.12
$ws.Range("I137").Value = 200787.2
$ws.Range("K137").Value = 602361.6000000001
$ws.Range("M137").Value = -597261.6000000001
$ws.Range("H139").Value = 45121.207
$ws.Range("I139").Value = 64231.875
$ws.Range("K139").Value = 192695.625
$ws.Range("M139").Value = -187555.625
$ws.Range("H140").Value = 250543.5
$ws.Range("I140").Value = 250543.5
$ws.Range("K140").Value = 751630.5
$ws.Range("M140").Value = -746450.5

# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H26").Value = 42069
$ws.Range("J26").Value = 42069
$ws.Range("L26").Value = 42069
$ws.Range("N26").Value = -42629
$ws.Range("H50").Value = 42069
$ws.Range("J50").Value = 42069
$ws.Range("L50").Value = 42069
$ws.Range("N50").Value = -43065
$ws.Range("H52").Value = 89991.2
$ws.Range("J52").Value = 89991.2
$ws.Range("L52").Value = 89991.2
$ws.Range("N52").Value = -90509.2
$ws.Range("H107").Value = 728427.0600000001
$ws.Range("I107").Value = 1334200
$ws.Range("K107").Value = 1334200
$ws.Range("M107").Value = -1332280

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("M4").ClearContents()
$ws.Range("H12").Value = 7925.3335
$ws.Range("J12").Value = 7999.5
$ws.Range("L12").Value = 7999.5
$ws.Range("N12").Value = -8339.5
$ws.Range("H28").Value = 0
$ws.Range("I28").Value = 0
$ws.Range("K28").Value = 0
$ws.Range("M28").ClearContents()
$ws.Range("H37").Value = 0
$ws.Range("I37").Value = 0
$ws.Range("K37").Value = 0
$ws.Range("M37").ClearContents()
$ws.Range("H93").Value = 3384.25
$ws.Range("I93").Value = 3232.375
$ws.Range("J93").Value = 3536.125
$ws.Range("K93").Value = 3232.375
$ws.Range("L93").Value = 3536.125
$ws.Range("M93").Value = -1984.375
$ws.Range("N93").Value = -6032.125
$ws.Range("H132").Value = 16673471
$ws.Range("I132").Value = 26321976
$ws.Range("K132").Value = 78965928
$ws.Range("M132").Value = -78963398
$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()
$ws.Range("H136").Value = 7464.304
$ws.Range("I136").Value = 6468.4
$ws.Range("J136").Value = 8230.385
$ws.Range("K136").Value = 19405.2
$ws.Range("L136").Value = 24691.155
$ws.Range("M136").Value = -16855.2
$ws.Range("N136").Value = -29791.155

# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 22001.5
$ws.Range("I15").Value = 22001.5
$ws.Range("K15").Value = 22001.5
$ws.Range("M15").Value = -21713.5
$ws.Range("H132").Value = 14295438
$ws.Range("I132").Value = 15154118
$ws.Range("K132").Value = 45462354
$ws.Range("M132").Value = -45459824
$ws.Range("H133").Value = 170000
$ws.Range("J133").Value = 170000
$ws.Range("L133").Value = 170000
$ws.Range("N133").Value = -180120
$ws.Range("H136").Value = 32293780
$ws.Range("I136").Value = 50001320
$ws.Range("K136").Value = 150003960
$ws.Range("M136").Value = -150001410
